$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Spankulator")

$ws.Range("A4").Value = "RV2 "
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = "20K"
$ws.Range("D4").Value = "https://www.amazon.com/dp/B015RK03LI?psc=1&ref=ppx_yo2ov_dt_b_product_details"
